$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D target cells to be treated as text so values like "256.34" are not
# reinterpreted as numbers (they must remain text, matching the original inlineStr cells).
$dRows = 2,3,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,31,32,33,34,35,36,37,38,39,40,43,44,45,46,48,49,50,51
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "98.272.09"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.418.27"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "256.34"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "665.88"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").Value = "0.437"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("D9").Value = "1.06"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "3.418.71"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("E12").Value = "  +3.66%  "
$ws.Range("D13").Value = "42.78"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "6.48"
$ws.Range("E14").Value = "  +16.94%  "
$ws.Range("D15").Value = "97.969.06"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "0.0000268"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "4.064.20"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "9.17"
$ws.Range("E18").Value = "  +24.46%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.419.77"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("B20").Value = "Stellar"
$ws.Range("C20").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D20").Value = "0.583"
$ws.Range("E20").Value = "  +34.00%  "
$ws.Range("D21").Value = "17.81"
$ws.Range("E21").Value = "  +6.24%  "
$ws.Range("D22").Value = "11.03"
$ws.Range("E22").Value = "  +8.55%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.45"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "511.40"
$ws.Range("E24").Value = "  -4.35%  "
$ws.Range("D25").Value = "0.0000207"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "6.50"
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("D27").Value = "101.14"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "12.99"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "3.603.27"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "11.73"
$ws.Range("E31").Value = "  +7.67%  "
$ws.Range("D32").Value = "0.200"
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "0.584"
$ws.Range("E34").Value = "  +11.29%  "
$ws.Range("D35").Value = "2.40"
$ws.Range("E35").Value = "  +16.72%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "30.10"
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  +16.67%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "7.98"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").Value = "538.46"
$ws.Range("E40").Value = "  +4.92%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "0.886"
$ws.Range("E43").Value = "  +9.87%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "24.71"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "9.15"
$ws.Range("E45").Value = "  +18.83%  "
$ws.Range("D46").Value = "5.88"
$ws.Range("E46").Value = "  +19.46%  "
$ws.Range("E47").Value = "  +7.93%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +17.38%  "
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "3.70"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "3.33"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "54.69"
$ws.Range("E51").Value = "  +10.62%  "

# Restore default (unstyled) cell style on the D cells we touched, so formatting
# matches the original workbook (no explicit style index).
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
